# Applies the Swahili (Kenya) translation edits described by the diff.
#
# We target each paragraph by its (stable) 0-based position in
# $d.Paragraphs and overwrite its Range.Text in one shot. Doing the
# replacement through the paragraph's own Range (rather than a
# Content-wide Find/Replace, or a Range manufactured via $d.Range(s,e))
# is what makes the engine keep each run's existing
# xml:space="preserve" handling intact instead of re-minting the <w:t>
# without it. A couple of cells contain a manual line break
# (<w:br/>) between two lines of text in the same run -- we reproduce
# that with a literal vertical-tab (`v, Word's in-memory line-break
# character) inside the replacement string.

$d = $word.ActiveDocument

$replacements = @{
    1  = "Kichwa cha Video"
    2  = "Tatizo la Uwanja wa Ndege"
    5  = "Mada"
    6  = "Jiometri"
    9  = "Malengo"
    10 = "Pata wazo angavu la tatizo la kupunguza, tambua jinsi ya kutekeleza kwa vitendo matatizo ya kupunguza."
    13 = "Urefu"
    17 = "Mahali pa Kambi"
    21 = "Wawezeshaji"
    25 = "N. ya wanafunzi"
    29 = "Tarehe"
    33 = "Rasilimali"
    34 = "inahitajika"
    35 = "Pini (3 kila kikundi), kamba (1/kikundi), pete ya chuma (si lazima lakini inafaa kuzuia msuguano 1/kikundi), kadibodi nene au sehemu ya mbao inayoweza kutupwa (1/kikundi)"
    38 = "Maandalizi"
    39 = "Weka alama 3 kwenye kuni"
    44 = "Muda wa video"
    45 = "Mwezeshaji anafanya nini"
    46 = "Wanachofanya wanafunzi"
    49 = "Utangulizi Mkuu wa Video ya VMC"
    53 = "Utangulizi wa Video"
    57 = "Kitendawili"
    61 = "Utangulizi wa jaribio la kwanza"
    64 = "SITISHA VIDEO`vKutafuta suluhu"
    65 = "Kusaidia mchakato, kuchochea mawazo"
    66 = "Jaribu kupata mpangilio wa kamba ili upunguzaji wa urefu wa kamba ulingane na kupunguza jumla ya urefu wa barabara"
    69 = "Suluhisho "
    72 = "SITISHA VIDEO`vJiometri"
    73 = "Kusaidia mchakato, kuchochea mawazo"
    74 = "Jaribu kujua ni mali gani ya kijiometri ambayo hatua mpya ina uhusiano na 3 ya kuanzia."
    77 = "Inaonyesha pembe 120°"
    83 = "suluhisho"
    85 = "Pointi ,F, inayopatikana kama inavyoonyeshwa kwenye video, inaitwa Fermat Point."
    86 = "Kuna njia kadhaa tofauti zinazowezekana za kamba ambazo zinaweza kutumika kupata uhakika F."
    87 = "Ona kwamba matumizi ya pete sio lazima kabisa, lakini inasaidia kupunguza msuguano (adui wa asili wa uzoefu huu)."
    88 = "Mara tu uhakika unapopatikana (kabla ya kutazama suluhu katika video) wanafunzi wanaweza kuulizwa kutafuta pembe kwa kutambua kwamba kila pembe ina mshikamano na zinaunda 360° zote kwa pamoja."
    89 = "Muundo wa kijiometri ambao unaweza kutumika kupata F ni pamoja na kujenga pembetatu zilizo sawa kwenye kando ya pembetatu asilia na kuunganisha alama tofauti:"
    91 = "Ujenzi huu unaweza kuigwa kwenye ubao wa mbao ili kuthibitisha kwamba njia hizo mbili zitaongoza kwenye hatua sawa."
}

$i = 0
foreach ($p in $d.Paragraphs) {
    if ($replacements.ContainsKey($i)) {
        $p.Range.Text = $replacements[$i]
    }
    $i++
}

# --- Document default language: sw-TZ -> sw-KE -------------------------
$d.Styles("Normal").LanguageID = "sw-KE"
